# Generate Report for Handoff
#
# Updates the "Latest Handoff Datetime" column (D) for the row
# corresponding to 0dbd8e8f-e224-47f0-a54f-fa8691f05d4a.md on both the
# zh-cn and de-de localization-status sheets, reflecting a newly
# generated handoff report/timestamp for that file.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D6").Value = "2016-03-09 22:38:26"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D6").Value = "2016-03-09 22:38:32"
